# Updates cryptos list values (Price and Volume(1h) columns) per latest data refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.069.32"
$ws.Range("E2").Value = "  -1.92%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.668.06"
$ws.Range("E3").Value = "  -1.57%  "
$ws.Range("E4").Value = "  -0.13%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "216.95"
$ws.Range("E5").Value = "  -1.07%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5110"
$ws.Range("E6").Value = "  +0.67%  "
$ws.Range("E7").Value = "  -0.14%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2653"
$ws.Range("E8").Value = "  +0.20%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06407"
$ws.Range("E9").Value = "  +2.08%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "21.92"
$ws.Range("E10").Value = "  -0.82%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07430"
$ws.Range("E11").Value = "  +0.80%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.690.84"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.501"
$ws.Range("E13").Value = "  -0.17%  "
$ws.Range("E14").Value = "  +0.44%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.000008552"
$ws.Range("E15").Value = "  +2.04%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "64.35"
$ws.Range("E16").Value = "  -1.70%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "25.974.63"
$ws.Range("E17").Value = "  -2.48%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "4.945"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.005"
$ws.Range("E19").Value = "  -0.07%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "10.77"
$ws.Range("E20").Value = "  -2.03%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "190.66"
$ws.Range("E21").Value = "  +2.68%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.226"
$ws.Range("E22").Value = "  -0.68%  "
$ws.Range("E23").Value = "  -0.13%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "145.33"
$ws.Range("E24").Value = "  +0.38%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "7.618"
$ws.Range("E25").Value = "  +1.30%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.1200"
$ws.Range("E26").Value = "  +4.08%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "15.63"
$ws.Range("E27").Value = "  -0.03%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.06745"
$ws.Range("E28").Value = "  +19.37%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.317"
$ws.Range("E29").Value = "  -1.49%  "
$ws.Range("E30").Value = "  -1.52%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.540"
$ws.Range("E31").Value = "  +0.63%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.523"
$ws.Range("E32").Value = "  +0.83%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.648"
$ws.Range("E33").Value = "  +0.25%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.018"
$ws.Range("E34").Value = "  -0.13%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.6103"
$ws.Range("E35").Value = "  +1.11%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.367"
$ws.Range("E36").Value = "  +0.08%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.717"
$ws.Range("E37").Value = "  +1.31%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "6.237"
$ws.Range("E38").Value = "  +6.58%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01603"
$ws.Range("E39").Value = "  -0.38%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.086.98"
$ws.Range("E40").Value = "  -0.86%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.8693"
$ws.Range("E41").Value = "  +0.82%  "
$ws.Range("E42").Value = "  +0.63%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "100.79"
$ws.Range("E43").Value = "  +1.07%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.817.90"
$ws.Range("E44").Value = "  -1.87%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00000000113"
$ws.Range("E45").Value = "  +1.00%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "56.40"
$ws.Range("E46").Value = "  -0.44%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.005"
$ws.Range("E47").Value = "  -0.04%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "8.095"
$ws.Range("E48").Value = "  -0.63%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.05236"
$ws.Range("E49").Value = "  -0.11%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.4286"
$ws.Range("E50").Value = "  -1.00%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "6.013"
$ws.Range("E51").Value = "  +4.11%  "
